# Weekly update: insert 3 new price rows (week of 44595) at the top of the
# "Sandia" price block, pushing the existing rows (378-386) down to 381-389.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 378 (existing row 378 and below shift down by 3).
$ws.Rows(378).Resize(3).Insert()

# Shared/static field values for every data row in this sheet.
$marketId   = 9
$market     = "Vega Central Mapocho de Santiago"
$region     = "Metropolitana"
$codreg     = 13
$categoryId = 100112028
$category   = "Sandia"
$variety    = "Sin especificar"
$unit       = "`$/unidad"
$kgOrUnits  = 1
$classif    = "Hortaliza"
$newDate    = 44595

# New row 378: Extra
$r = 378
$ws.Cells.Item($r, 1).Value  = $marketId
$ws.Cells.Item($r, 2).Value  = $market
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $newDate
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $categoryId
$ws.Cells.Item($r, 7).Value  = $category
$ws.Cells.Item($r, 8).Value  = $variety
$ws.Cells.Item($r, 9).Value  = "Extra"
$ws.Cells.Item($r, 10).Value = 300
$ws.Cells.Item($r, 11).Value = 3500
$ws.Cells.Item($r, 12).Value = 3500
$ws.Cells.Item($r, 13).Value = 3500
$ws.Cells.Item($r, 14).Value = $unit
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 3500
$ws.Cells.Item($r, 17).Value = $kgOrUnits
$ws.Cells.Item($r, 18).Value = $classif

# New row 379: Primera
$r = 379
$ws.Cells.Item($r, 1).Value  = $marketId
$ws.Cells.Item($r, 2).Value  = $market
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $newDate
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $categoryId
$ws.Cells.Item($r, 7).Value  = $category
$ws.Cells.Item($r, 8).Value  = $variety
$ws.Cells.Item($r, 9).Value  = "Primera"
$ws.Cells.Item($r, 10).Value = 450
$ws.Cells.Item($r, 11).Value = 3000
$ws.Cells.Item($r, 12).Value = 3000
$ws.Cells.Item($r, 13).Value = 3000
$ws.Cells.Item($r, 14).Value = $unit
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 3000
$ws.Cells.Item($r, 17).Value = $kgOrUnits
$ws.Cells.Item($r, 18).Value = $classif

# New row 380: Segunda
$r = 380
$ws.Cells.Item($r, 1).Value  = $marketId
$ws.Cells.Item($r, 2).Value  = $market
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $newDate
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $categoryId
$ws.Cells.Item($r, 7).Value  = $category
$ws.Cells.Item($r, 8).Value  = $variety
$ws.Cells.Item($r, 9).Value  = "Segunda"
$ws.Cells.Item($r, 10).Value = 220
$ws.Cells.Item($r, 11).Value = 2500
$ws.Cells.Item($r, 12).Value = 2500
$ws.Cells.Item($r, 13).Value = 2500
$ws.Cells.Item($r, 14).Value = $unit
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 2500
$ws.Cells.Item($r, 17).Value = $kgOrUnits
$ws.Cells.Item($r, 18).Value = $classif
